# Error Calculations and Plots
# This script reproduces the target edit:
#  1. Two rows are removed from the missing-data sheet ("RM 232" and "SC 92"),
#     which shifts all subsequent rows up and shrinks the used range from
#     A1:F35 down to A1:F33.
#  2. A number of individual cells flip between a present numeric value and a
#     "missing" (blank / inline-string) placeholder, reflecting an updated
#     imputation/removal pattern for the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the two rows that disappear entirely ---
# Delete the lower-indexed row first is fine as long as we delete by the
# ORIGINAL row numbers in the right order (higher row number first keeps the
# other row index valid).
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"

# --- Step 2: apply the individual cell content changes (post row-shift row numbers) ---
$ws.Range("C3").Value = 11.2
$ws.Range("D4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("D11").Value = -15.5
$ws.Range("D12").Value = -14.1
$ws.Range("F12").Value = ""
$ws.Range("F13").Value = 17.1
$ws.Range("F14").Value = 17.76
$ws.Range("D15").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = 17.78
$ws.Range("D18").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("D31").Value = -13.7
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("F32").Value = 17.39

Write-Output $ws.UsedRange.Address()
